$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 35, shifting existing rows 35-132 down to 36-133.
$ws.Rows(35).Insert()

# Populate the newly inserted row 35 with the new record.
$ws.Range("A35").Value = 9
$ws.Range("B35").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44811
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = 100112022
$ws.Range("G35").Value = "Arveja Verde"
$ws.Range("H35").Value = "Perfection"
$ws.Range("I35").Value = "Primera"
$ws.Range("J35").Value = 25
$ws.Range("K35").Value = 36000
$ws.Range("L35").Value = 38000
$ws.Range("M35").Value = 37200
$ws.Range("N35").Value = "$/malla 25 kilos"
$ws.Range("O35").Value = "Provincia de Huasco"
$ws.Range("P35").Value = 1488
$ws.Range("Q35").Value = 25
$ws.Range("R35").Value = "Hortaliza"
